$wb = $excel.ActiveWorkbook

$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsCreateAccount = $wb.Worksheets.Item("CreateAccount")

# The sign-up email changed (23/03/2018 -> 18/04/2018 test address). The same
# text is shared between SignIn!A2 and CreateAccount!D2, so update both cells
# to the new address.
$wsSignIn.Range("A2").Value = "testjaga18042018@gmail.com"
$wsCreateAccount.Range("D2").Value = "testjaga18042018@gmail.com"

# The email on the CreateAccount sheet should also be a clickable mailto link
# (fixes the drop-down/email link issue). Re-apply the original "Hyperlink"
# cell style afterwards since Add() otherwise pushes its own duplicate style.
$wsCreateAccount.Hyperlinks.Add($wsCreateAccount.Range("D2"), "mailto:testjaga18042018@gmail.com") | Out-Null
$wsCreateAccount.Range("D2").Style = $wsCreateAccount.Range("L2").Style

# Update each sheet's remembered selection.
$wsSignIn.Range("A2").Select() | Out-Null
$wsCreateAccount.Range("F3").Select() | Out-Null

# CreateAccount is now the active (foreground) sheet/tab.
$wsCreateAccount.Activate()

$wb.Save()
